$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D: copy the existing header formatting from C1
# (bold font, thin border all around, centered/top aligned) onto D1, then
# set its text.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "success"

# Determine which row has the max "count" (column C) among the data rows (2..43)
$maxRow = 2
$maxVal = $ws.Cells.Item(2, 3).Value2
for ($r = 3; $r -le 43; $r++) {
    $v = $ws.Cells.Item($r, 3).Value2
    if ($v -gt $maxVal) {
        $maxVal = $v
        $maxRow = $r
    }
}

# Fill column D ("success") for each data row: "1" for the max-count row, "0" otherwise
# (leading apostrophe forces Excel to store the value as literal text, matching
# the text/shared-string "0"/"1" entries used for the success flag)
for ($r = 2; $r -le 43; $r++) {
    if ($r -eq $maxRow) {
        $ws.Cells.Item($r, 4).Formula = "'1"
    } else {
        $ws.Cells.Item($r, 4).Formula = "'0"
    }
}
